# Update scaling mappings to proper format
# Rewrites the "year" sheet of the Edgar scaling mapping workbook: splits the
# combined multi-year comment rows for rou/mkd into one row per year, moves the
# start/end scaling year values into the F/G columns for idn and adds a new
# svk row, reordering rows to match the new canonical layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("year")

# Clear out the old data rows (keep header row 1) before rewriting.
$ws.Range("A2:H9").ClearContents()

$header = @("iso","scaling_sector","pre_ext_year","post_ext_year","select_scaling_year","start_scaling_year","end_scaling_year","Comment")
for ($c = 1; $c -le $header.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $header[$c - 1]
}

$rows = @(
    @("idn","1A2","NA","NA","NA",1990,2010,"Eliminate jump in CO emissions present in EDGAR data"),
    @("svk","all","NA","NA","NA",2000,2010,"Don't calibrate to flat 1990s flat portion that is not in EMEP expert data"),
    @("rou","all","NA","NA",1980,"NA","NA","Reduce jumps in emissions"),
    @("rou","all","NA","NA",1991,"NA","NA","Reduce jumps in emissions"),
    @("rou","all","NA","NA",2000,"NA","NA","Reduce jumps in emissions"),
    @("rou","all","NA","NA",2010,"NA","NA","Reduce jumps in emissions"),
    @("mkd","all","NA","NA",1980,"NA","NA","Reduce jumps in emissions"),
    @("mkd","all","NA","NA",1992,"NA","NA","Reduce jumps in emissions")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

$ws.Range("A10:XFD11").Select() | Out-Null
